# Add a default header to the document showing "Questionnaire 7" so the
# questionnaire number stays visible after printing.

$d = $word.ActiveDocument

# Primary (default) header of the document's single section.
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)   # wdHeaderFooterPrimary = 1

# Insert the header text (InsertAfter avoids materializing the
# even-page/first-page header & footer variants that a plain
# "Range.Text = ..." assignment would create).
$hdr.Range.InsertAfter("Questionnaire 7")

# Paragraph formatting: built-in "Header" style, centered.
$hdr.Range.Style = "Header"
$hdr.Range.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

# Character formatting on the text itself (exclude the trailing
# paragraph mark so no extra rPr ends up on the pPr).
$txt = $hdr.Range.Duplicate
$txt.End = $txt.End - 1
$txt.Font.Name = "Arial"
$txt.Font.Size = 12
